# Add habitat_typenumber (F) and mfd_hab1/2/3 (N/O/P) values for P01_2 rows 2-45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column F to text format since several habitat type numbers are purely numeric
# (e.g. "7230", "6403") and would otherwise be auto-converted to numbers.
$ws.Range("F2:F45").NumberFormat = "@"

$ws.Range("F2").Value = "91D0"
$ws.Range("N2").Value = "Forests"
$ws.Range("O2").Value = "Temperate forests"
$ws.Range("P2").Value = "Bog woodland"
$ws.Range("F3").Value = "91E0"
$ws.Range("N3").Value = "Forests"
$ws.Range("O3").Value = "Temperate forests"
$ws.Range("P3").Value = "Alluvial woodland"
$ws.Range("F4").Value = "91E0"
$ws.Range("N4").Value = "Forests"
$ws.Range("O4").Value = "Temperate forests"
$ws.Range("P4").Value = "Alluvial woodland"
$ws.Range("F5").Value = "9990"
$ws.Range("N5").Value = "Forests"
$ws.Range("O5").Value = "Forest (non-habitattype)"
$ws.Range("P5").Value = "Deciduous trees (løvtræer)"
$ws.Range("F6").Value = "9920"
$ws.Range("N6").Value = "Forests"
$ws.Range("O6").Value = "Forest (non-habitattype)"
$ws.Range("P6").Value = "Willow"
$ws.Range("F7").Value = "91D0"
$ws.Range("N7").Value = "Forests"
$ws.Range("O7").Value = "Temperate forests"
$ws.Range("P7").Value = "Bog woodland"
$ws.Range("F8").Value = "9920"
$ws.Range("N8").Value = "Forests"
$ws.Range("O8").Value = "Forest (non-habitattype)"
$ws.Range("P8").Value = "Willow"
$ws.Range("F9").Value = "9920"
$ws.Range("N9").Value = "Forests"
$ws.Range("O9").Value = "Forest (non-habitattype)"
$ws.Range("P9").Value = "Willow"
$ws.Range("F10").Value = "7230"
$ws.Range("N10").Value = "Bogs, mires and fens"
$ws.Range("O10").Value = "Calcareous fens"
$ws.Range("P10").Value = "Alkaline fens"
$ws.Range("F11").Value = "7230"
$ws.Range("N11").Value = "Bogs, mires and fens"
$ws.Range("O11").Value = "Calcareous fens"
$ws.Range("P11").Value = "Alkaline fens"
$ws.Range("F12").Value = "6403"
$ws.Range("N12").Value = "Grassland formations"
$ws.Range("O12").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P12").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F13").Value = "7230"
$ws.Range("N13").Value = "Bogs, mires and fens"
$ws.Range("O13").Value = "Calcareous fens"
$ws.Range("P13").Value = "Alkaline fens"
$ws.Range("F14").Value = "6403"
$ws.Range("N14").Value = "Grassland formations"
$ws.Range("O14").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P14").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F15").Value = "7230"
$ws.Range("N15").Value = "Bogs, mires and fens"
$ws.Range("O15").Value = "Calcareous fens"
$ws.Range("P15").Value = "Alkaline fens"
$ws.Range("F16").Value = "6403"
$ws.Range("N16").Value = "Grassland formations"
$ws.Range("O16").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P16").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F17").Value = "6403"
$ws.Range("N17").Value = "Grassland formations"
$ws.Range("O17").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P17").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F18").Value = "91D0"
$ws.Range("N18").Value = "Forests"
$ws.Range("O18").Value = "Temperate forests"
$ws.Range("P18").Value = "Bog woodland"
$ws.Range("F19").Value = "6403"
$ws.Range("N19").Value = "Grassland formations"
$ws.Range("O19").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P19").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F20").Value = "7007"
$ws.Range("N20").Value = "Bogs, mires and fens"
$ws.Range("O20").Value = "Fen wetland (non-habitat type)"
$ws.Range("P20").Value = "Wet fens"
$ws.Range("F21").Value = "7007"
$ws.Range("N21").Value = "Bogs, mires and fens"
$ws.Range("O21").Value = "Fen wetland (non-habitat type)"
$ws.Range("P21").Value = "Wet fens"
$ws.Range("F22").Value = "6402"
$ws.Range("N22").Value = "Grassland formations"
$ws.Range("O22").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P22").Value = "Natural meadow (6410 subtype)"
$ws.Range("F23").Value = "6402"
$ws.Range("N23").Value = "Grassland formations"
$ws.Range("O23").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P23").Value = "Natural meadow (6410 subtype)"
$ws.Range("F24").Value = "6403"
$ws.Range("N24").Value = "Grassland formations"
$ws.Range("O24").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P24").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F25").Value = "6402"
$ws.Range("N25").Value = "Grassland formations"
$ws.Range("O25").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P25").Value = "Natural meadow (6410 subtype)"
$ws.Range("F26").Value = "7230"
$ws.Range("N26").Value = "Bogs, mires and fens"
$ws.Range("O26").Value = "Calcareous fens"
$ws.Range("P26").Value = "Alkaline fens"
$ws.Range("F27").Value = "7230"
$ws.Range("N27").Value = "Bogs, mires and fens"
$ws.Range("O27").Value = "Calcareous fens"
$ws.Range("P27").Value = "Alkaline fens"
$ws.Range("F28").Value = "91D0"
$ws.Range("N28").Value = "Forests"
$ws.Range("O28").Value = "Temperate forests"
$ws.Range("P28").Value = "Bog woodland"
$ws.Range("F29").Value = "6403"
$ws.Range("N29").Value = "Grassland formations"
$ws.Range("O29").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P29").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F30").Value = "6402"
$ws.Range("N30").Value = "Grassland formations"
$ws.Range("O30").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P30").Value = "Natural meadow (6410 subtype)"
$ws.Range("F31").Value = "91D0"
$ws.Range("N31").Value = "Forests"
$ws.Range("O31").Value = "Temperate forests"
$ws.Range("P31").Value = "Bog woodland"
$ws.Range("F32").Value = "6402"
$ws.Range("N32").Value = "Grassland formations"
$ws.Range("O32").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P32").Value = "Natural meadow (6410 subtype)"
$ws.Range("F33").Value = "6403"
$ws.Range("N33").Value = "Grassland formations"
$ws.Range("O33").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P33").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F34").Value = "7140"
$ws.Range("N34").Value = "Bogs, mires and fens"
$ws.Range("O34").Value = "Sphagnum acid bogs"
$ws.Range("P34").Value = "Quaking bogs"
$ws.Range("F35").Value = "9920"
$ws.Range("N35").Value = "Forests"
$ws.Range("O35").Value = "Forest (non-habitattype)"
$ws.Range("P35").Value = "Willow"
$ws.Range("F36").Value = "6410"
$ws.Range("N36").Value = "Grassland formations"
$ws.Range("O36").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P36").Value = "Molinia meadows"
$ws.Range("F37").Value = "6402"
$ws.Range("N37").Value = "Grassland formations"
$ws.Range("O37").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P37").Value = "Natural meadow (6410 subtype)"
$ws.Range("F38").Value = "7900"
$ws.Range("N38").Value = "Bogs, mires and fens"
$ws.Range("O38").Value = "§3 mire"
$ws.Range("P38").Value = "§3 mire"
$ws.Range("F39").Value = "7900"
$ws.Range("N39").Value = "Bogs, mires and fens"
$ws.Range("O39").Value = "§3 mire"
$ws.Range("P39").Value = "§3 mire"
$ws.Range("F40").Value = "6403"
$ws.Range("N40").Value = "Grassland formations"
$ws.Range("O40").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P40").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F41").Value = "6403"
$ws.Range("N41").Value = "Grassland formations"
$ws.Range("O41").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P41").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F42").Value = "6402"
$ws.Range("N42").Value = "Grassland formations"
$ws.Range("O42").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P42").Value = "Natural meadow (6410 subtype)"
$ws.Range("F43").Value = "6403"
$ws.Range("N43").Value = "Grassland formations"
$ws.Range("O43").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P43").Value = "Agricultural meadow (6430 subtype)"
$ws.Range("F44").Value = "6402"
$ws.Range("N44").Value = "Grassland formations"
$ws.Range("O44").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P44").Value = "Natural meadow (6410 subtype)"
$ws.Range("F45").Value = "6402"
$ws.Range("N45").Value = "Grassland formations"
$ws.Range("O45").Value = "Semi-natural tall-herb humid meadows"
$ws.Range("P45").Value = "Natural meadow (6410 subtype)"
